$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 48
$ws.Range("H48").Value = 3063.8
$ws.Range("J48").Value = 3329.75
$ws.Range("L48").Value = 9989.25
$ws.Range("N48").Value = -10573.25

# Row 56
$ws.Range("H56").Value = 3063.8
$ws.Range("J56").Value = 3329.75
$ws.Range("L56").Value = 9989.25
$ws.Range("N56").Value = -11057.25

# Row 93
$ws.Range("H93").Value = 76360.39999999999
$ws.Range("J93").Value = 76360.39999999999
$ws.Range("L93").Value = 76360.39999999999
$ws.Range("N93").Value = -81352.39999999999

# Row 128
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960

# Row 132
$ws.Range("H132").Value = 919.4808
$ws.Range("I132").Value = 829.3333
$ws.Range("K132").Value = 2487.9999
$ws.Range("M132").Value = 42.0001000000002

# Row 138
$ws.Range("H138").Value = 1844.84
$ws.Range("I138").Value = 1353.7241
$ws.Range("J138").Value = 2045.4366
$ws.Range("K138").Value = 4061.1723
$ws.Range("L138").Value = 6136.3098
$ws.Range("M138").Value = 1078.8277
$ws.Range("N138").Value = -16416.3098

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6646.35
$ws.Range("I32").Value = 4043.814
$ws.Range("J32").Value = 22633.357
$ws.Range("K32").Value = 4043.814
$ws.Range("L32").Value = 22633.357
$ws.Range("M32").Value = -3756.814
$ws.Range("N32").Value = -23207.357

# Row 61
$ws.Range("H61").Value = 31537.777
$ws.Range("I61").Value = 37053.273
$ws.Range("J61").Value = 7269.6
$ws.Range("K61").Value = 37053.273
$ws.Range("L61").Value = 7269.6
$ws.Range("M61").Value = -36841.273
$ws.Range("N61").Value = -7693.6

# Row 102
$ws.Range("H102").Value = 994.2857
$ws.Range("I102").Value = 994.2857
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 994.2857
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 627.7143
$ws.Range("N102").ClearContents()

# Row 132
$ws.Range("H132").Value = 1904.8392
$ws.Range("I132").Value = 1435.3704
$ws.Range("J132").Value = 2341.9312
$ws.Range("K132").Value = 4306.1112
$ws.Range("L132").Value = 7025.7936
$ws.Range("M132").Value = -1776.1112
$ws.Range("N132").Value = -12085.7936

# Row 136
$ws.Range("H136").Value = 31537.777
$ws.Range("I136").Value = 37053.273
$ws.Range("J136").Value = 7269.6
$ws.Range("K136").Value = 111159.819
$ws.Range("L136").Value = 21808.8
$ws.Range("M136").Value = -108609.819
$ws.Range("N136").Value = -26908.8

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1915.3125
$ws.Range("I99").Value = 1895.2307
$ws.Range("K99").Value = 1895.2307
$ws.Range("M99").Value = -397.2307000000001

# Row 105
$ws.Range("H105").Value = 2424.8386
$ws.Range("I105").Value = 2272.3333
$ws.Range("K105").Value = 2272.3333
$ws.Range("M105").Value = -525.3332999999998

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2929.3125
$ws.Range("I31").Value = 2170.1667
$ws.Range("K31").Value = 2170.1667
$ws.Range("M31").Value = -1875.1667

# Row 34
$ws.Range("H34").Value = 2929.3125
$ws.Range("I34").Value = 2170.1667
$ws.Range("K34").Value = 2170.1667
$ws.Range("M34").Value = -1968.1667

# Row 86
$ws.Range("H86").Value = 2306.7693
$ws.Range("I86").Value = 2061.25
$ws.Range("J86").Value = 2699.6
$ws.Range("K86").Value = 2061.25
$ws.Range("L86").Value = 2699.6
$ws.Range("M86").Value = -938.25
$ws.Range("N86").Value = -4945.6

# Row 89
$ws.Range("H89").Value = 2306.7693
$ws.Range("I89").Value = 2061.25
$ws.Range("J89").Value = 2699.6
$ws.Range("K89").Value = 10306.25
$ws.Range("L89").Value = 13498
$ws.Range("M89").Value = -4690.25
$ws.Range("N89").Value = -24730

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1066.6666
$ws.Range("I68").Value = 700
$ws.Range("K68").Value = 2100
$ws.Range("M68").Value = -1289

# Row 71
$ws.Range("H71").Value = 1066.6666
$ws.Range("I71").Value = 700
$ws.Range("K71").Value = 6300
$ws.Range("M71").Value = -2244

# Row 129
$ws.Range("H129").Value = 28245.518
$ws.Range("J129").Value = 30426.76
$ws.Range("L129").Value = 91280.28
$ws.Range("N129").Value = -101280.28

# Row 131
$ws.Range("H131").Value = 23899.742
$ws.Range("I131").Value = 433
$ws.Range("J131").Value = 26099.75
$ws.Range("K131").Value = 1299
$ws.Range("L131").Value = 78299.25
$ws.Range("M131").Value = 3741
$ws.Range("N131").Value = -88379.25

# Row 141
$ws.Range("H141").Value = 3141.6667
$ws.Range("I141").Value = 3286.7
$ws.Range("J141").Value = 2416.5
$ws.Range("K141").Value = 9860.099999999999
$ws.Range("L141").Value = 7249.5
$ws.Range("M141").Value = -4680.099999999999
$ws.Range("N141").Value = -17609.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4927.5
$ws.Range("I70").Value = 4600
$ws.Range("J70").Value = 5255
$ws.Range("K70").Value = 4600
$ws.Range("L70").Value = 5255
$ws.Range("M70").Value = -4330
$ws.Range("N70").Value = -5795

# Row 73
$ws.Range("H73").Value = 4927.5
$ws.Range("I73").Value = 4600
$ws.Range("J73").Value = 5255
$ws.Range("K73").Value = 4600
$ws.Range("L73").Value = 5255
$ws.Range("M73").Value = -3664
$ws.Range("N73").Value = -7127

# Row 102
$ws.Range("H102").Value = 2089
$ws.Range("I102").Value = 2273.1428
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 2273.1428
$ws.Range("L102").Value = 800
$ws.Range("M102").Value = -651.1428000000001
$ws.Range("N102").Value = -4044

# Row 122
$ws.Range("H122").Value = 1434
$ws.Range("J122").Value = 1406
$ws.Range("L122").Value = 4218
$ws.Range("N122").Value = -9118

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 1568.3
$ws.Range("I100").Value = 1076
$ws.Range("K100").Value = 1076
$ws.Range("M100").Value = -535

# Row 122
$ws.Range("H122").Value = 2517.2273
$ws.Range("I122").Value = 2367.4666
$ws.Range("K122").Value = 7102.399800000001
$ws.Range("M122").Value = -4652.399800000001

# Row 132
$ws.Range("H132").Value = 2395.4736
$ws.Range("I132").Value = 1434.08
$ws.Range("J132").Value = 4244.3076
$ws.Range("K132").Value = 4302.24
$ws.Range("L132").Value = 12732.9228
$ws.Range("M132").Value = -1772.24
$ws.Range("N132").Value = -17792.9228

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 617.61536
$ws.Range("I113").Value = 457.18182
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1371.54546
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = 798.45454
$ws.Range("N113").Value = -8840

# Row 122
$ws.Range("H122").Value = 69664.5
$ws.Range("I122").Value = 91330.44500000001
$ws.Range("K122").Value = 273991.335
$ws.Range("M122").Value = -271541.335

# Row 128
$ws.Range("H128").Value = 32500
$ws.Range("J128").Value = 32500
$ws.Range("L128").Value = 32500
$ws.Range("N128").Value = -42460
